$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.142.91'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.636.80'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.60'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.516'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.14%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.97'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0846'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.867.16'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.658.11'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.12'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.540'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.61'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.158.44'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0739'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '216.66'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.41%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.53%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.10'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.09%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.38'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.299.83'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.78%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.45'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.544'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.38%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.852'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.65%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.806'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.89%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.777.01'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.60'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.27'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0107'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.61'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0955'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.44%  '
